$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move/update the active selection to I16
$ws.Range("I16").Select()

# Adjust column widths (no columns are inserted - letters B:G keep their data).
# Values are chosen so the engine's ColumnWidth->stored-width pixel rounding
# lands on the same target pixel widths as the authored file
# (stored width = round(ColumnWidth*6 + 5)/6 for this engine).
$ws.Columns("B").ColumnWidth = 2.5
$ws.Columns("C").ColumnWidth = 14
$ws.Columns("D").ColumnWidth = 17
$ws.Columns("E").ColumnWidth = 35.5
$ws.Columns("F").ColumnWidth = 15.333333333333334
